$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9; this shifts existing rows 9-58 down to 10-59,
# preserving all of their values (a new weekly record is being added).
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new weekly record.
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C9").Value = "Arica y Parinacota"
$ws.Range("D9").Value = 45061
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100108
$ws.Range("H9").Value = "Tropicales y subtropicales"
$ws.Range("I9").Value = 100108001
$ws.Range("J9").Value = "Guayaba"
$ws.Range("K9").Value = "Sin especificar"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 40
$ws.Range("N9").Value = 600
$ws.Range("O9").Value = 700
$ws.Range("P9").Value = 638
$ws.Range("Q9").Value = "$/kilo (en caja de 10 kilos )"
$ws.Range("R9").Value = "Región de Arica y Parinacota"
$ws.Range("S9").Value = 638
$ws.Range("T9").Value = 1
